$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.146.13"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "2.267.89"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "498.35"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.01"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.79%  "
$ws.Range("E9").Value = "  +0.32%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.334"
$ws.Range("E11").Value = "  +3.35%  "
$ws.Range("E12").Value = "  +2.52%  "
$ws.Range("D13").Value = "2.670.16"
$ws.Range("E13").Value = "  -0.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.64"
$ws.Range("E14").Value = "  +3.94%  "
$ws.Range("D15").Value = "54.124.73"
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "2.267.55"
$ws.Range("E17").Value = "  -1.73%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.20"
$ws.Range("E18").Value = "  +1.65%  "
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "302.86"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("E23").Value = "  -3.13%  "
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("E26").Value = "  +2.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "170.86"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0692"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.60"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("E32").Value = "  +0.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.74"
$ws.Range("E33").Value = "  +0.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.937"
$ws.Range("E35").Value = "  +9.02%  "
$ws.Range("E36").Value = "  -1.24%  "
$ws.Range("E37").Value = "  +1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.372"
$ws.Range("E38").Value = "  -0.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.39"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("E40").Value = "  +0.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.79"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "124.49"
$ws.Range("E42").Value = "  -3.35%  "
$ws.Range("E43").Value = "  +1.81%  "
$ws.Range("E44").Value = "  -0.44%  "
$ws.Range("E45").Value = "  -0.87%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "238.64"
$ws.Range("E46").Value = "  -1.50%  "
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  +0.82%  "
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("E50").Value = "  -1.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.64"
